$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.409.51"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.65%  "
$ws.Range("D3").Value = "'2.586.84"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.07%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'590.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.89%  "
$ws.Range("D6").Value = "'150.81"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.49%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  -0.45%  "
$ws.Range("E9").Value = "  +0.69%  "
$ws.Range("E10").Value = "  +1.94%  "
$ws.Range("D11").Value = "'0.384"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.67%  "
$ws.Range("E12").Value = "  -0.36%  "
$ws.Range("D13").Value = "'27.60"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.56%  "
$ws.Range("D14").Value = "'3.049.32"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.12%  "
$ws.Range("D15").Value = "'63.194.97"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.73%  "
$ws.Range("E16").Value = "  +5.86%  "
$ws.Range("D17").Value = "'2.587.30"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.14%  "
$ws.Range("D18").Value = "'12.24"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.90%  "
$ws.Range("E19").Value = "  +3.69%  "
$ws.Range("D20").Value = "'345.57"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.47%  "
$ws.Range("E21").Value = "  -0.60%  "
$ws.Range("E22").Value = "  -0.04%  "
$ws.Range("D23").Value = "'67.47"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.72%  "
$ws.Range("E24").Value = "  +1.52%  "
$ws.Range("E25").Value = "  +0.44%  "
$ws.Range("D26").Value = "'1.67"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.72%  "
$ws.Range("D27").Value = "'561.15"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.69%  "
$ws.Range("D28").Value = "'8.07"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.47%  "
$ws.Range("E29").Value = "  +0.96%  "
$ws.Range("E30").Value = "  +0.27%  "
$ws.Range("D31").Value = "'2.03"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.24%  "
$ws.Range("D32").Value = "'0.0₃0854"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.54%  "
$ws.Range("E33").Value = "  -0.64%  "
$ws.Range("E34").Value = "  -0.85%  "
$ws.Range("D35").Value = "'166.74"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.63%  "
$ws.Range("E36").Value = "  +2.18%  "
$ws.Range("E37").Value = "  +0.13%  "
$ws.Range("E38").Value = "  +1.91%  "
$ws.Range("D39").Value = "'1.92"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.65%  "
$ws.Range("D41").Value = "'166.28"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.57%  "
$ws.Range("D42").Value = "'39.59"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.44%  "
$ws.Range("D43").Value = "'4.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.48%  "
$ws.Range("D44").Value = "'22.82"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.96%  "
$ws.Range("E45").Value = "  +3.12%  "
$ws.Range("D46").Value = "'2.12"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.89%  "
$ws.Range("E47").Value = "  +0.20%  "
$ws.Range("E48").Value = "  +2.50%  "
$ws.Range("E49").Value = "  +0.74%  "
$ws.Range("D50").Value = "'19.26"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.21%  "
$ws.Range("D51").Value = "'0.0₆0233"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +18.05%  "
